# Update worksheet values with new TPM-derived numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> column -> new value
$updates = @{
    2 = @{
        'E' = 2
        'F' = 0.6666666666666666
        'G' = 0.3971766666666667
        'H' = 1.19153
        'M' = 0.481166
        'N' = 1.443498
        'O' = 0.1587222767546368
        'P' = 0.1690961013073894
        'Q' = 0.1911079079933333
        'R' = 1.71997117194
        'S' = 0.1587222767546368
        'T' = 0.1690961013073894
    }
    3 = @{
        'E' = 2
        'F' = 0.6666666666666666
        'G' = 0.3971766666666667
        'H' = 1.19153
        'M' = 1.046115666666667
        'N' = 3.138347
        'O' = 0.3450822800489395
        'P' = 0.3676362850864647
        'Q' = 0.4154927334344445
        'R' = 3.73943460091
        'S' = 0.3450822800489395
        'T' = 0.3676362850864647
    }
    4 = @{
        'E' = 2
        'F' = 0.6666666666666666
        'G' = 0.3971766666666667
        'H' = 1.19153
        'M' = 0.6239283333333333
        'N' = 1.871785
        'O' = 0.2058153019922285
        'P' = 0.2192670485069268
        'Q' = 0.2478097756722222
        'R' = 2.23028798105
        'S' = 0.2058153019922285
        'T' = 0.2192670485069268
    }
    5 = @{
        'E' = 2
        'F' = 0.6666666666666666
        'G' = 0.3971766666666667
        'H' = 1.19153
        'M' = 0.5579350000000001
        'N' = 1.11587
        'O' = 0.1840460744963241
        'P' = 0.130716680290431
        'Q' = 0.2215987635166667
        'R' = 1.3295925811
        'S' = 0.1840460744963241
        'T' = 0.130716680290431
    }
    6 = @{
        'E' = 2
        'F' = 0.6666666666666666
        'G' = 0.3971766666666667
        'H' = 1.19153
        'M' = 0.3223513333333334
        'N' = 0.9670540000000001
        'O' = 0.1063340667078711
        'P' = 0.1132838848087882
        'Q' = 0.1280304280688889
        'R' = 1.15227385262
        'S' = 0.1063340667078711
        'T' = 0.1132838848087882
    }
}

foreach ($rowNum in $updates.Keys) {
    $cols = $updates[$rowNum]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$rowNum").Value = $cols[$col]
    }
}
